$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''41.947.57'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.02%  '

$ws.Range('D3').Value = '''2.228.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.90%  '

$ws.Range('E4').Value = '  +0.13%  '

$ws.Range('D5').Value = '''231.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.79%  '

$ws.Range('D6').Value = '''0.625'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.55%  '

$ws.Range('D7').Value = '''61.06'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.64%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').Value = '''0.402'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.19%  '

$ws.Range('E10').Value = '  +0.96%  '

$ws.Range('E11').Value = '  +5.06%  '

$ws.Range('E12').Value = '  -0.18%  '

$ws.Range('D13').Value = '''2.559.44'
$ws.Range('D13').Style = 'Normal'

$ws.Range('D14').Value = '''15.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.19%  '

$ws.Range('E15').Value = '  +0.63%  '

$ws.Range('D16').Value = '''0.799'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.76%  '

$ws.Range('D17').Value = '''5.57'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.19%  '

$ws.Range('D18').Value = '''2.249.80'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.83%  '

$ws.Range('D19').Value = '''41.851.46'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.74%  '

$ws.Range('D20').Value = '''72.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.42%  '

$ws.Range('D21').Value = '''0.0₃0890'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.20%  '

$ws.Range('D22').Value = '''6.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.60%  '

$ws.Range('D23').Value = '''249.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.58%  '

$ws.Range('E24').Value = '  -0.03%  '

$ws.Range('E25').Value = '  +1.63%  '

$ws.Range('D26').Value = '''2.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.11%  '

$ws.Range('D27').Value = '''9.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.56%  '

$ws.Range('D28').Value = '''0.142'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.45%  '

$ws.Range('D29').Value = '''167.43'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.70%  '

$ws.Range('D30').Value = '''19.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.50%  '

$ws.Range('E31').Value = '  -2.13%  '

$ws.Range('E32').Value = '  -1.70%  '

$ws.Range('E33').Value = '  -0.09%  '

$ws.Range('E34').Value = '  +5.58%  '

$ws.Range('E35').Value = '  +3.36%  '

$ws.Range('E36').Value = '  +1.97%  '

$ws.Range('E37').Value = '  -4.47%  '

$ws.Range('E38').Value = '  -2.85%  '

$ws.Range('E39').Value = '  -1.29%  '

$ws.Range('E40').Value = '  +31.14%  '

$ws.Range('E41').Value = '  +0.01%  '

$ws.Range('E42').Value = '  +5.69%  '

$ws.Range('D43').Value = '''4.79'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.37%  '

$ws.Range('D44').Value = '''8.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +8.71%  '

$ws.Range('E45').Value = '  +7.24%  '

$ws.Range('D46').Value = '''99.09'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.89%  '

$ws.Range('E47').Value = '  +0.67%  '

$ws.Range('D48').Value = '''1.470.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.76%  '

$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').Value = '''2.81'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.27%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '''16.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.04%  '

$ws.Range('E51').Value = '  -1.14%  '
